$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J column values for rows 2-10
$values = @(
    @{Row=2;  I=1; J=1},
    @{Row=3;  I=1; J=5},
    @{Row=4;  I=1; J=2},
    @{Row=5;  I=6; J=6},
    @{Row=6;  I=1; J=4},
    @{Row=7;  I=1; J=4},
    @{Row=8;  I=1; J=3},
    @{Row=9;  I=5; J=6},
    @{Row=10; I=3; J=3}
)

foreach ($entry in $values) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}
